$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($addr, $val) {
    # Force the cell to hold a TEXT value even when the string looks like a
    # number (e.g. "593.73"), matching the inlineStr/shared-string cells the
    # source workbook already uses for the Price / Volume(1h) columns.
    # Temporarily apply a Text number format so Excel doesn't auto-coerce
    # the assignment into a float, then clear the formatting back off again
    # so the cell's style index returns to its original (default) value.
    $cell = $ws.Range($addr)
    $cell.NumberFormat = "@"
    $cell.Value = $val
    $cell.ClearFormats()
}

# --- Rows 42 and 43 swap ranking positions (Maker now above EnergySwap) ---
# Before: Row42=EnergySwap(26.14,-1.28%), Row43=Maker(2.799.35,+0.57%)
# After:  Row42=Maker(2.813.24,+0.94%),   Row43=EnergySwap(26.14,-1.48%)
$ws.Range("B42").Value = "Maker"
$ws.Range("C42").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
Set-TextValue "D42" "2.813.24"
$ws.Range("E42").Value = "  +0.94%  "

$ws.Range("B43").Value = "EnergySwap"
$ws.Range("C43").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
Set-TextValue "D43" "26.14"
$ws.Range("E43").Value = "  -1.48%  "

# --- Price / Volume(1h) updates for all other rows ---
Set-TextValue "D2" "66.858.02"
$ws.Range("E2").Value = "  +0.35%  "

Set-TextValue "D3" "3.490.79"
$ws.Range("E3").Value = "  -0.06%  "

$ws.Range("E4").Value = "  -0.02%  "

Set-TextValue "D5" "593.73"
$ws.Range("E5").Value = "  +0.46%  "

Set-TextValue "D6" "171.76"
$ws.Range("E6").Value = "  +1.40%  "

$ws.Range("E7").Value = "  +0.01%  "

$ws.Range("E8").Value = "  -0.76%  "

$ws.Range("E9").Value = "  +3.36%  "

$ws.Range("E11").Value = "  -2.11%  "

Set-TextValue "D12" "4.093.95"
$ws.Range("E12").Value = "  -0.03%  "

$ws.Range("E13").Value = "  -0.51%  "

Set-TextValue "D14" "28.89"
$ws.Range("E14").Value = "  +2.36%  "

Set-TextValue "D15" "66.842.79"
$ws.Range("E15").Value = "  +0.32%  "

$ws.Range("E16").Value = "  -0.47%  "

Set-TextValue "D17" "3.513.62"
$ws.Range("E17").Value = "  +0.97%  "

Set-TextValue "D18" "6.26"
$ws.Range("E18").Value = "  -1.05%  "

Set-TextValue "D19" "14.02"
$ws.Range("E19").Value = "  -0.85%  "

Set-TextValue "D20" "392.78"
$ws.Range("E20").Value = "  +0.14%  "

Set-TextValue "D21" "7.95"
$ws.Range("E21").Value = "  -0.39%  "

Set-TextValue "D22" "72.89"
$ws.Range("E22").Value = "  -0.65%  "

$ws.Range("E23").Value = "  +0.03%  "

$ws.Range("E24").Value = "  -0.48%  "

$ws.Range("E25").Value = "  -1.53%  "

Set-TextValue "D26" "10.14"
$ws.Range("E26").Value = "  -0.25%  "

$ws.Range("E27").Value = "  -0.80%  "

$ws.Range("E28").Value = "  -0.31%  "

Set-TextValue "D29" "6.16"
$ws.Range("E29").Value = "  -2.97%  "

Set-TextValue "D30" "1.43"
$ws.Range("E30").Value = "  -3.30%  "

$ws.Range("E31").Value = "  -0.70%  "

Set-TextValue "D32" "23.69"
$ws.Range("E32").Value = "  +0.29%  "

Set-TextValue "D33" "7.33"
$ws.Range("E33").Value = "  -1.23%  "

$ws.Range("E34").Value = "  -0.61%  "

Set-TextValue "D35" "162.88"
$ws.Range("E35").Value = "  +0.29%  "

Set-TextValue "D36" "0.877"
$ws.Range("E36").Value = "  -0.83%  "

$ws.Range("E37").Value = "  -1.93%  "

Set-TextValue "D38" "6.93"
$ws.Range("E38").Value = "  +2.39%  "

Set-TextValue "D39" "4.65"
$ws.Range("E39").Value = "  -0.80%  "

Set-TextValue "D40" "0.0739"
$ws.Range("E40").Value = "  -1.00%  "

Set-TextValue "D41" "27.14"
$ws.Range("E41").Value = "  -0.47%  "

Set-TextValue "D44" "42.81"
$ws.Range("E44").Value = "  -0.93%  "

Set-TextValue "D45" "2.55"
$ws.Range("E45").Value = "  +1.72%  "

Set-TextValue "D46" "0.0302"
$ws.Range("E46").Value = "  -3.82%  "

Set-TextValue "D47" "338.37"
$ws.Range("E47").Value = "  -3.32%  "

Set-TextValue "D48" "34.75"
$ws.Range("E48").Value = "  +2.13%  "

Set-TextValue "D49" "1.08"
$ws.Range("E49").Value = "  -2.41%  "

$ws.Range("E50").Value = "  -1.31%  "

Set-TextValue "D51" "6.41"
$ws.Range("E51").Value = "  -2.33%  "
